$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32
$ws.Range("A32").Value = 112212309
$ws.Range("B32").Value = 90814
$ws.Range("E32").Value = 4364
$ws.Range("F32").Value = "Dropptaggsvamp"
$ws.Range("G32").Value = "Hydnellum ferrugineum"
$ws.Range("H32").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I32").Value = "'2"
$ws.Range("I32").Style = "Normal"
$ws.Range("J32").Value = "fruktkroppar"
$ws.Range("P32").Value = "Simsbodarna SO, Dlr"
$ws.Range("Q32").Value = 515492
$ws.Range("R32").Value = 6704591
$ws.Range("Z32").Value = "12:08"
$ws.Range("AB32").Value = "12:08"

# Row 33
$ws.Range("A33").Value = 112212369
$ws.Range("B33").Value = 90814
$ws.Range("D33").Value = "LC"
$ws.Range("E33").Value = 4364
$ws.Range("F33").Value = "Dropptaggsvamp"
$ws.Range("G33").Value = "Hydnellum ferrugineum"
$ws.Range("H33").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I33").Value = "'10"
$ws.Range("I33").Style = "Normal"
$ws.Range("Q33").Value = 515497
$ws.Range("R33").Value = 6704555
$ws.Range("S33").Value = 2
$ws.Range("Z33").Value = "12:12"
$ws.Range("AB33").Value = "12:12"

# Row 34
$ws.Range("A34").Value = 112212094
$ws.Range("P34").Value = "Simsbodarna SO, Dlr"
$ws.Range("Q34").Value = 515450
$ws.Range("R34").Value = 6704585
$ws.Range("S34").Value = 2
$ws.Range("Z34").Value = "11:59"
$ws.Range("AB34").Value = "11:59"
$ws.Range("AC34").Value = "Ca. Färska fk."

# Row 35
$ws.Range("A35").Value = 112212286
$ws.Range("B35").Value = 90837
$ws.Range("D35").Value = "NT"
$ws.Range("E35").Value = 5966
$ws.Range("F35").Value = "Motaggsvamp"
$ws.Range("G35").Value = "Sarcodon squamosus"
$ws.Range("H35").Value = "(Schaeff.) Quél."
$ws.Range("I35").Value = "'3"
$ws.Range("I35").Style = "Normal"
$ws.Range("Q35").Value = 515476
$ws.Range("R35").Value = 6704606
$ws.Range("S35").Value = 6
$ws.Range("Z35").Value = "12:08"
$ws.Range("AB35").Value = "12:08"

# Row 36
$ws.Range("A36").Value = 112211876
$ws.Range("B36").Value = 90837
$ws.Range("D36").Value = "NT"
$ws.Range("E36").Value = 5966
$ws.Range("F36").Value = "Motaggsvamp"
$ws.Range("G36").Value = "Sarcodon squamosus"
$ws.Range("H36").Value = "(Schaeff.) Quél."
$ws.Range("P36").Value = "SO Simsbodarna, Dlr"
$ws.Range("Q36").Value = 515371
$ws.Range("R36").Value = 6704633
$ws.Range("Z36").Value = "11:43"
$ws.Range("AB36").Value = "11:43"

# Row 37
$ws.Range("A37").Value = 112212788
$ws.Range("Q37").Value = 515101
$ws.Range("R37").Value = 6704641
$ws.Range("Z37").Value = "12:37"
$ws.Range("AB37").Value = "12:37"

# Row 38
$ws.Range("A38").Value = 112212237
$ws.Range("I38").Value = "'2"
$ws.Range("I38").Style = "Normal"
$ws.Range("P38").Value = "Simsbodarna SO, Dlr"
$ws.Range("Q38").Value = 515492
$ws.Range("R38").Value = 6704591
$ws.Range("Z38").Value = "11:59"
$ws.Range("AB38").Value = "11:59"
$ws.Range("AC38").Value = "Stora fina ex."

# Row 39
$ws.Range("A39").Value = 112211016
$ws.Range("I39").Value = "'1"
$ws.Range("I39").Style = "Normal"
$ws.Range("P39").Value = "S Simsbodarna, Dlr"
$ws.Range("Q39").Value = 515070
$ws.Range("R39").Value = 6704834
$ws.Range("Z39").Value = "10:55"
$ws.Range("AB39").Value = "10:55"

# Row 40
$ws.Range("A40").Value = 112211348
$ws.Range("B40").Value = 89517
$ws.Range("E40").Value = 5447
$ws.Range("F40").Value = "Vedticka"
$ws.Range("G40").Value = "Fuscoporia viticola"
$ws.Range("H40").Value = "(Schwein.) Murrill"
$ws.Range("I40").Value = ""
$ws.Range("J40").Value = ""
$ws.Range("P40").Value = "Simsbodarna S, Dlr"
$ws.Range("Q40").Value = 515173
$ws.Range("R40").Value = 6704768
$ws.Range("S40").Value = 1
$ws.Range("Z40").Value = "11:15"
$ws.Range("AB40").Value = "11:15"
$ws.Range("AC40").Value = ""

# Row 41
$ws.Range("A41").Value = 112211929
$ws.Range("B41").Value = 90814
$ws.Range("D41").Value = "LC"
$ws.Range("E41").Value = 4364
$ws.Range("F41").Value = "Dropptaggsvamp"
$ws.Range("G41").Value = "Hydnellum ferrugineum"
$ws.Range("H41").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I41").Value = "'10"
$ws.Range("I41").Style = "Normal"
$ws.Range("P41").Value = "SO Simsbodarna, Dlr"
$ws.Range("Q41").Value = 515370
$ws.Range("R41").Value = 6704610
$ws.Range("S41").Value = 5
$ws.Range("Z41").Value = "11:45"
$ws.Range("AB41").Value = "11:45"
$ws.Range("AC41").Value = "G:a fk."
